$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B..J, rows 2..21 (row 7 update - Premier League corner stats)
$data = @{
    2  = @(7, 4, 41, 10.25, 3, 3, 9, 50, 7.142857142857143)
    3  = @(7, 4, 22, 5.5, 3, 2.333333333333333, 7, 29, 4.142857142857143)
    4  = @(7, 3, 20, 6.666666666666667, 4, 7, 28, 48, 6.857142857142857)
    5  = @(7, 4, 23, 5.75, 3, 3.666666666666667, 11, 34, 4.857142857142857)
    6  = @(7, 4, 26, 6.5, 3, 4.666666666666667, 14, 40, 5.714285714285714)
    7  = @(7, 4, 27, 6.75, 3, 4.333333333333333, 13, 40, 5.714285714285714)
    8  = @(7, 4, 38, 9.5, 3, 6.666666666666667, 20, 58, 8.285714285714286)
    9  = @(7, 4, 15, 3.75, 3, 6.333333333333333, 19, 34, 4.857142857142857)
    10 = @(7, 4, 14, 3.5, 3, 2.666666666666667, 8, 22, 3.142857142857143)
    11 = @(7, 3, 16, 5.333333333333333, 4, 5.5, 22, 38, 5.428571428571429)
    12 = @(7, 3, 20, 6.666666666666667, 4, 4.5, 18, 38, 5.428571428571429)
    13 = @(7, 4, 9, 2.25, 3, 2.333333333333333, 7, 16, 2.285714285714286)
    14 = @(7, 3, 19, 6.333333333333333, 4, 7.5, 30, 49, 7)
    15 = @(7, 3, 17, 5.666666666666667, 4, 6.25, 25, 42, 6)
    16 = @(7, 3, 15, 5, 4, 6.5, 26, 41, 5.857142857142857)
    17 = @(7, 3, 15, 5, 4, 4, 16, 31, 4.428571428571429)
    18 = @(7, 3, 6, 2, 4, 6.25, 25, 31, 4.428571428571429)
    19 = @(7, 3, 28, 9.333333333333334, 4, 8.75, 35, 63, 9)
    20 = @(7, 4, 19, 4.75, 3, 3.666666666666667, 11, 30, 4.285714285714286)
    21 = @(7, 3, 11, 3.666666666666667, 4, 3.25, 13, 24, 3.428571428571428)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 2 + $i  # column B = 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
